$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 10 from 45183 to 45184
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
